$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking values
# (e.g. "581.34") are preserved exactly as strings rather than
# being re-interpreted as floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.104.91"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.593.64"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "581.34"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "191.27"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("D8").Value = "3.590.93"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "0.183"
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").Value = "0.665"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "55.94"
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("D13").Value = "0.0000307"
$ws.Range("E13").Value = "  +6.94%  "
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "4.177.38"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "19.99"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("D17").Value = "3.597.66"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "70.131.26"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "12.70"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "482.31"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "19.26"
$ws.Range("E23").Value = "  +10.47%  "
$ws.Range("E24").Value = "  -6.18%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "95.42"
$ws.Range("E26").Value = "  +5.69%  "
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "9.43"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "32.19"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "7.67"
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("D33").Value = "12.24"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "66.76"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "584.19"
$ws.Range("E35").Value = "  -6.31%  "
$ws.Range("D36").Value = "39.00"
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "0.0₃0805"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "0.398"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Value = "3.32"
$ws.Range("E40").Value = "  +23.63%  "
$ws.Range("E41").Value = "  -3.84%  "
$ws.Range("D42").Value = "3.225.69"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  -6.37%  "
$ws.Range("E44").Value = "  +7.08%  "
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "0.0451"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").Value = "9.48"
$ws.Range("E47").Value = "  +4.57%  "
$ws.Range("D48").Value = "3.32"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "0.138"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").Value = "  -4.77%  "

# Restore the default (unstyled) cell style on column D so no
# stray style index is left behind on these cells.
$ws.Range("D2:D51").Style = "Normal"
